# Apply updated crypto price/volume data to sheet (matches author commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.318.74"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = "'3.420.74"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'576.38"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.32%  '
$ws.Range('D6').Value = "'128.35"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.72%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'0.479"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.58%  '
$ws.Range('D9').Value = "'7.51"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.29%  '
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').Value = "'4.004.12"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.20%  '
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('E14').Value = '  -2.87%  '
$ws.Range('D15').Value = "'3.421.18"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('D16').Value = "'63.384.15"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.51%  '
$ws.Range('D17').Value = "'25.06"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').Value = "'9.79"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').Value = "'5.64"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.72%  '
$ws.Range('D20').Value = "'13.15"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.49%  '
$ws.Range('D21').Value = "'382.41"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('D22').Value = "'0.560"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.94%  '
$ws.Range('D23').Value = "'3.558.94"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.11%  '
$ws.Range('D24').Value = "'73.58"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('E26').Value = '  -5.12%  '
$ws.Range('D27').Value = "'0.987"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('E28').Value = '  -3.34%  '
$ws.Range('D29').Value = "'6.99"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.81%  '
$ws.Range('D30').Value = "'7.87"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.94%  '
$ws.Range('D31').Value = "'0.152"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').Value = "'1.40"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.56%  '
$ws.Range('D33').Value = "'3.451.74"
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Value = "'22.59"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E36').Value = '  +0.68%  '
$ws.Range('E37').Value = '  -2.09%  '
$ws.Range('D38').Value = "'164.09"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('E39').Value = '  -2.88%  '
$ws.Range('D40').Value = "'0.0760"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.39%  '
$ws.Range('D41').Value = "'0.783"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.10%  '
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').Value = "'41.17"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.32%  '
$ws.Range('D44').Value = "'4.29"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('E45').Value = '  -4.02%  '
$ws.Range('D46').Value = "'23.11"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.68%  '
$ws.Range('E47').Value = '  -6.39%  '
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('E49').Value = '  -1.15%  '
$ws.Range('D50').Value = "'2.261.68"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.71%  '
$ws.Range('E51').Value = '  -2.94%  '
